$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers (losing literal formatting like trailing zeros).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated values / labels / links cell by cell.
$ws.Range('D2').Value = '61.737.71'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '3.384.76'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '574.36'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '138.40'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.384.03'
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').Value = '0.475'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = '7.49'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').Value = '0.125'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '0.392'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '3.950.98'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').Value = '0.0000176'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '26.11'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').Value = '3.371.66'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '61.807.18'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = '14.00'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '9.42'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '379.21'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '0.557'
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').Value = '3.512.94'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '0.0000127'
$ws.Range('E26').Value = '  +5.41%  '
$ws.Range('D27').Value = '71.53'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('D28').Value = '1.81'
$ws.Range('E28').Value = '  +10.42%  '
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').Value = '0.996'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('D31').Value = '0.164'
$ws.Range('E31').Value = '  +4.69%  '
$ws.Range('D32').Value = '8.30'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '23.65'
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').Value = '5.26'
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('D37').Value = '6.87'
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('D38').Value = '1.54'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '165.00'
$ws.Range('E39').Value = '  +2.74%  '
$ws.Range('D40').Value = '0.0774'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '1.74'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = '0.774'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('D45').Value = '41.65'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = '4.40'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').Value = '24.58'
$ws.Range('E47').Value = '  +6.36%  '
$ws.Range('D48').Value = '6.89'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('D49').Value = '23.10'
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = '2.358.78'
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('D51').Value = '0.0265'
$ws.Range('E51').Value = '  +0.91%  '
